## Fall 2015/Ideas.docx -- "Added Lectures, Updated Syllabus"
##
## 1. Marks "Adwords" (in the LP-formulations sentence) as a spell-check
##    exception by splitting the single run into three runs and wrapping
##    the middle one in <w:proofErr w:type="spellStart"/> ... spellEnd/>,
##    matching what Word does automatically whenever it re-spell-checks a
##    capitalized word it doesn't recognise.
## 2. Moves the (hidden) _GoBack bookmark that used to sit at the end of
##    that paragraph into its own, otherwise-empty paragraph.
## 3. Appends a new paragraph with the "For algorithms ..." lecture note.
##
## The host COM emulation only supports structural paragraph splits via a
## whole-story XML replace (Range.InsertParagraphAfter/Before here always
## operate at whole-paragraph granularity), so the three edits are applied
## together as a single Range.InsertXML over the full document content.

$d = $word.ActiveDocument

$body = $d.Content
$body.Find.ClearFormatting()
$body.Find.Execute(
    "For the LP formulations do the network revenue management problem and also include the Adwords problem.",
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

if (-not $body.Find.Found) {
    throw "Could not locate the target sentence to edit."
}

$newBodyXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p w14:paraId="7B43B603" w14:textId="77777777" w:rsidR="00B401B8" w:rsidRDefault="00A5676C">
            <w:r>
              <w:t xml:space="preserve">For the LP formulations do the network revenue management problem and also include the </w:t>
            </w:r>
            <w:proofErr w:type="spellStart"/>
            <w:r>
              <w:t>Adwords</w:t>
            </w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:r>
              <w:t xml:space="preserve"> problem.</w:t>
            </w:r>
          </w:p>
          <w:p w14:paraId="5A2F4E10" w14:textId="5A2F4E10" w:rsidR="00B401B8" w:rsidRDefault="00B401B8">
            <w:bookmarkStart w:id="0" w:name="_GoBack"/>
            <w:bookmarkEnd w:id="0"/>
          </w:p>
          <w:p w14:paraId="3C9D7B24" w14:textId="3C9D7B24" w:rsidR="00B401B8" w:rsidRDefault="00B401B8">
            <w:r>
              <w:t>For algorithms and hard problems do LP relaxations and greedy algorithms</w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$body.InsertXML($newBodyXml) | Out-Null

Write-Output ("Paragraphs.Count=" + $d.Paragraphs.Count)
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    Write-Output ("Para " + $i + ": [" + $p.Range.Text + "]")
}
